$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the board-template placeholder cell (write-date -> reg-date)
$ws.Range("D4").Value = '${listview.regdate}'

# Column A no longer carries an explicit (unused/default) cell style;
# clear its format without touching any cell content/dimension.
$ws.Range("A:A").ClearFormats()
$ws.Range("A1").ClearContents()
$ws.Range("A3").ClearContents()
$ws.Range("A4").ClearContents()

# Update the active selection to D4 (matches the saved view state)
$ws.Range("D4").Select()
